# Update the NTT Data address block on the cover page (first page) of the
# WBS Dictionary document.
#
#   "NTT Data Romania"            -> "NTT Data Romania S.A."
#   "Street Constanta 19-21 "     -> "19-21, Constanta Street,"
#   "Cluj Napoca City, 400158"    -> "400158 Cluj Napoca"

$d = $word.ActiveDocument

# wdFindContinue=1, wdReplaceOne=2 -- replace just the single occurrence of
# each string, keeping the run/formatting the text lives in untouched.
$null = $d.Content.Find.Execute(
    "NTT Data Romania", $false, $false, $false, $false, $false,
    $true, 1, $false, "NTT Data Romania S.A.", 2)

$null = $d.Content.Find.Execute(
    "Street Constanta 19-21 ", $false, $false, $false, $false, $false,
    $true, 1, $false, "19-21, Constanta Street,", 2)

$null = $d.Content.Find.Execute(
    "Cluj Napoca City, 400158", $false, $false, $false, $false, $false,
    $true, 1, $false, "400158 Cluj Napoca", 2)
